$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1548163333333333
$ws.Cells.Item(2, 8).Value = 0.464449
$ws.Cells.Item(2, 9).Value = 0.04476815357596578
$ws.Cells.Item(2, 10).Value = 0.04476815357596578
$ws.Cells.Item(2, 13).Value = 211.7369283333333
$ws.Cells.Item(2, 14).Value = 635.210785
$ws.Cells.Item(2, 15).Value = 0.471139807893958
$ws.Cells.Item(2, 16).Value = 0.471139807893958
$ws.Cells.Item(2, 17).Value = 32.78033487582945
$ws.Cells.Item(2, 18).Value = 295.023013882465
$ws.Cells.Item(2, 19).Value = 0.02109205927554773
$ws.Cells.Item(2, 20).Value = 0.02109205927554773
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1548163333333333
$ws.Cells.Item(3, 8).Value = 0.464449
$ws.Cells.Item(3, 9).Value = 0.04476815357596578
$ws.Cells.Item(3, 10).Value = 0.04476815357596578
$ws.Cells.Item(3, 13).Value = 76.35132866666667
$ws.Cells.Item(3, 14).Value = 229.053986
$ws.Cells.Item(3, 15).Value = 0.1698907724959131
$ws.Cells.Item(3, 16).Value = 0.1698907724959131
$ws.Cells.Item(3, 17).Value = 11.82043274930156
$ws.Cells.Item(3, 18).Value = 106.383894743714
$ws.Cells.Item(3, 19).Value = 0.0076056961942365
$ws.Cells.Item(3, 20).Value = 0.007605696194236499
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1548163333333333
$ws.Cells.Item(4, 8).Value = 0.464449
$ws.Cells.Item(4, 9).Value = 0.04476815357596578
$ws.Cells.Item(4, 10).Value = 0.04476815357596578
$ws.Cells.Item(4, 13).Value = 1.001605
$ws.Cells.Item(4, 14).Value = 3.004815
$ws.Cells.Item(4, 15).Value = 0.00222869005893356
$ws.Cells.Item(4, 16).Value = 0.00222869005893356
$ws.Cells.Item(4, 17).Value = 0.1550648135483333
$ws.Cells.Item(4, 18).Value = 1.395583321935
$ws.Cells.Item(4, 19).Value = 0.00009977433883156586
$ws.Cells.Item(4, 20).Value = 0.00009977433883156586
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.1548163333333333
$ws.Cells.Item(5, 8).Value = 0.464449
$ws.Cells.Item(5, 9).Value = 0.04476815357596578
$ws.Cells.Item(5, 10).Value = 0.04476815357596578
$ws.Cells.Item(5, 13).Value = 25.733869
$ws.Cells.Item(5, 14).Value = 77.201607
$ws.Cells.Item(5, 15).Value = 0.05726091425082595
$ws.Cells.Item(5, 16).Value = 0.05726091425082595
$ws.Cells.Item(5, 17).Value = 3.984023241060333
$ws.Cells.Item(5, 18).Value = 35.856209169543
$ws.Cells.Item(5, 19).Value = 0.002563465403081184
$ws.Cells.Item(5, 20).Value = 0.002563465403081184
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.1548163333333333
$ws.Cells.Item(6, 8).Value = 0.464449
$ws.Cells.Item(6, 9).Value = 0.04476815357596578
$ws.Cells.Item(6, 10).Value = 0.04476815357596578
$ws.Cells.Item(6, 13).Value = 134.590487
$ws.Cells.Item(6, 14).Value = 403.771461
$ws.Cells.Item(6, 15).Value = 0.2994798153003695
$ws.Cells.Item(6, 16).Value = 0.2994798153003695
$ws.Cells.Item(6, 17).Value = 20.83680569888767
$ws.Cells.Item(6, 18).Value = 187.531251289989
$ws.Cells.Item(6, 19).Value = 0.01340715836426881
$ws.Cells.Item(6, 20).Value = 0.01340715836426881
$ws.Cells.Item(7, 7).Value = 3.059667666666666
$ws.Cells.Item(7, 8).Value = 9.179003
$ws.Cells.Item(7, 9).Value = 0.8847624087429419
$ws.Cells.Item(7, 10).Value = 0.8847624087429419
$ws.Cells.Item(7, 13).Value = 211.7369283333333
$ws.Cells.Item(7, 14).Value = 635.210785
$ws.Cells.Item(7, 15).Value = 0.471139807893958
$ws.Cells.Item(7, 16).Value = 0.471139807893958
$ws.Cells.Item(7, 17).Value = 647.8446334608172
$ws.Cells.Item(7, 18).Value = 5830.601701147355
$ws.Cells.Item(7, 19).Value = 0.4168467912869452
$ws.Cells.Item(7, 20).Value = 0.4168467912869452
$ws.Cells.Item(8, 7).Value = 3.059667666666666
$ws.Cells.Item(8, 8).Value = 9.179003
$ws.Cells.Item(8, 9).Value = 0.8847624087429419
$ws.Cells.Item(8, 10).Value = 0.8847624087429419
$ws.Cells.Item(8, 13).Value = 76.35132866666667
$ws.Cells.Item(8, 14).Value = 229.053986
$ws.Cells.Item(8, 15).Value = 0.1698907724959131
$ws.Cells.Item(8, 16).Value = 0.1698907724959131
$ws.Cells.Item(8, 17).Value = 233.6096916284398
$ws.Cells.Item(8, 18).Value = 2102.487224655958
$ws.Cells.Item(8, 19).Value = 0.1503129690966832
$ws.Cells.Item(8, 20).Value = 0.1503129690966832
$ws.Cells.Item(9, 7).Value = 3.059667666666666
$ws.Cells.Item(9, 8).Value = 9.179003
$ws.Cells.Item(9, 9).Value = 0.8847624087429419
$ws.Cells.Item(9, 10).Value = 0.8847624087429419
$ws.Cells.Item(9, 13).Value = 1.001605
$ws.Cells.Item(9, 14).Value = 3.004815
$ws.Cells.Item(9, 15).Value = 0.00222869005893356
$ws.Cells.Item(9, 16).Value = 0.00222869005893356
$ws.Cells.Item(9, 17).Value = 3.064578433271666
$ws.Cells.Item(9, 18).Value = 27.581205899445
$ws.Cells.Item(9, 19).Value = 0.001971861184883506
$ws.Cells.Item(9, 20).Value = 0.001971861184883506
$ws.Cells.Item(10, 7).Value = 3.059667666666666
$ws.Cells.Item(10, 8).Value = 9.179003
$ws.Cells.Item(10, 9).Value = 0.8847624087429419
$ws.Cells.Item(10, 10).Value = 0.8847624087429419
$ws.Cells.Item(10, 13).Value = 25.733869
$ws.Cells.Item(10, 14).Value = 77.201607
$ws.Cells.Item(10, 15).Value = 0.05726091425082595
$ws.Cells.Item(10, 16).Value = 0.05726091425082595
$ws.Cells.Item(10, 17).Value = 78.73708691753566
$ws.Cells.Item(10, 18).Value = 708.6337822578209
$ws.Cells.Item(10, 19).Value = 0.05066230441938382
$ws.Cells.Item(10, 20).Value = 0.05066230441938382
$ws.Cells.Item(11, 7).Value = 3.059667666666666
$ws.Cells.Item(11, 8).Value = 9.179003
$ws.Cells.Item(11, 9).Value = 0.8847624087429419
$ws.Cells.Item(11, 10).Value = 0.8847624087429419
$ws.Cells.Item(11, 13).Value = 134.590487
$ws.Cells.Item(11, 14).Value = 403.771461
$ws.Cells.Item(11, 15).Value = 0.2994798153003695
$ws.Cells.Item(11, 16).Value = 0.2994798153003695
$ws.Cells.Item(11, 17).Value = 411.8021613148203
$ws.Cells.Item(11, 18).Value = 3706.219451833383
$ws.Cells.Item(11, 19).Value = 0.2649684827550463
$ws.Cells.Item(11, 20).Value = 0.2649684827550463
$ws.Cells.Item(12, 7).Value = 0.243696
$ws.Cells.Item(12, 8).Value = 0.731088
$ws.Cells.Item(12, 9).Value = 0.07046943768109237
$ws.Cells.Item(12, 10).Value = 0.07046943768109236
$ws.Cells.Item(12, 13).Value = 211.7369283333333
$ws.Cells.Item(12, 14).Value = 635.210785
$ws.Cells.Item(12, 15).Value = 0.471139807893958
$ws.Cells.Item(12, 16).Value = 0.471139807893958
$ws.Cells.Item(12, 17).Value = 51.59944248712
$ws.Cells.Item(12, 18).Value = 464.39498238408
$ws.Cells.Item(12, 19).Value = 0.0332009573314651
$ws.Cells.Item(12, 20).Value = 0.03320095733146509
$ws.Cells.Item(13, 7).Value = 0.243696
$ws.Cells.Item(13, 8).Value = 0.731088
$ws.Cells.Item(13, 9).Value = 0.07046943768109237
$ws.Cells.Item(13, 10).Value = 0.07046943768109236
$ws.Cells.Item(13, 13).Value = 76.35132866666667
$ws.Cells.Item(13, 14).Value = 229.053986
$ws.Cells.Item(13, 15).Value = 0.1698907724959131
$ws.Cells.Item(13, 16).Value = 0.1698907724959131
$ws.Cells.Item(13, 17).Value = 18.606513390752
$ws.Cells.Item(13, 18).Value = 167.458620516768
$ws.Cells.Item(13, 19).Value = 0.01197210720499339
$ws.Cells.Item(13, 20).Value = 0.01197210720499338
$ws.Cells.Item(14, 7).Value = 0.243696
$ws.Cells.Item(14, 8).Value = 0.731088
$ws.Cells.Item(14, 9).Value = 0.07046943768109237
$ws.Cells.Item(14, 10).Value = 0.07046943768109236
$ws.Cells.Item(14, 13).Value = 1.001605
$ws.Cells.Item(14, 14).Value = 3.004815
$ws.Cells.Item(14, 15).Value = 0.00222869005893356
$ws.Cells.Item(14, 16).Value = 0.00222869005893356
$ws.Cells.Item(14, 17).Value = 0.24408713208
$ws.Cells.Item(14, 18).Value = 2.19678418872
$ws.Cells.Item(14, 19).Value = 0.0001570545352184886
$ws.Cells.Item(14, 20).Value = 0.0001570545352184886
$ws.Cells.Item(15, 7).Value = 0.243696
$ws.Cells.Item(15, 8).Value = 0.731088
$ws.Cells.Item(15, 9).Value = 0.07046943768109237
$ws.Cells.Item(15, 10).Value = 0.07046943768109236
$ws.Cells.Item(15, 13).Value = 25.733869
$ws.Cells.Item(15, 14).Value = 77.201607
$ws.Cells.Item(15, 15).Value = 0.05726091425082595
$ws.Cells.Item(15, 16).Value = 0.05726091425082595
$ws.Cells.Item(15, 17).Value = 6.271240939824
$ws.Cells.Item(15, 18).Value = 56.44116845841599
$ws.Cells.Item(15, 19).Value = 0.004035144428360954
$ws.Cells.Item(15, 20).Value = 0.004035144428360953
$ws.Cells.Item(16, 7).Value = 0.243696
$ws.Cells.Item(16, 8).Value = 0.731088
$ws.Cells.Item(16, 9).Value = 0.07046943768109237
$ws.Cells.Item(16, 10).Value = 0.07046943768109236
$ws.Cells.Item(16, 13).Value = 134.590487
$ws.Cells.Item(16, 14).Value = 403.771461
$ws.Cells.Item(16, 15).Value = 0.2994798153003695
$ws.Cells.Item(16, 16).Value = 0.2994798153003695
$ws.Cells.Item(16, 17).Value = 32.799163319952
$ws.Cells.Item(16, 18).Value = 295.192469879568
$ws.Cells.Item(16, 19).Value = 0.02110417418105445
$ws.Cells.Item(16, 20).Value = 0.02110417418105444
